$wb = $excel.ActiveWorkbook

# --- Contacts sheet (sheet1): add a new column E with sheet-name values ---
$contacts = $wb.Worksheets.Item("Contacts")
$contacts.Range("E2").Value = "sheet1"
$contacts.Range("E3").Value = "Sheet2"
$contacts.Range("E4").Value = "Sheet 3"

# Update selection to match the new active cell on the Contacts sheet
$contacts.Activate()
$contacts.Range("E4").Select()

# --- Templates sheet (sheet2): just move the selection ---
$templates = $wb.Worksheets.Item("Templates")
$templates.Activate()
$templates.Range("C3").Select()

# Re-activate Contacts sheet as the tab shown (tabSelected stays on sheet1)
$contacts.Activate()
